# Add a new "loginUser" worksheet right after the existing "getUser" sheet,
# and make it the active/selected sheet (matching activeTab="1" in the workbook).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "loginUser"

# Header row
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("C1").Value = "statusCode"

# Row 2: valid login test case (A2 becomes a mailto hyperlink further below)
$newSheet.Range("A2").Value = "test@test.com"
$newSheet.Range("B2").Value = "test"
$newSheet.Range("C2").Value = 200

# Row 3: numeric username test case
$newSheet.Range("A3").Value = 12234
$newSheet.Range("B3").Value = "test"
$newSheet.Range("C3").Value = 200

# Rows are populated out of numeric order below (5, 6, then 4) so that the
# resulting shared-strings table is built up in the same sequence as the
# target workbook (Password, test@test.com, test, "", " ", SD21@@!#D).

# Row 5: empty-string credentials test case
$newSheet.Range("A5").Value = """"""
$newSheet.Range("B5").Value = """"""
$newSheet.Range("C5").Value = 200

# Row 6: whitespace credentials test case
$newSheet.Range("A6").Value = " "
$newSheet.Range("B6").Value = " "
$newSheet.Range("C6").Value = 200

# Row 4: bad/invalid password test case (no Password value)
$newSheet.Range("A4").Value = "SD21@@!#D"
$newSheet.Range("C4").Value = 200

# Turn A2 into a mailto hyperlink and restore the built-in "Hyperlink" cell
# style (Hyperlinks.Add applies its own style slot, so reassign explicitly).
$hlink = $newSheet.Hyperlinks.Add($newSheet.Range("A2"), "mailto:test@test.com")
$newSheet.Range("A2").Style = "Hyperlink"

# Column A width to fit the email/username values
$newSheet.Columns.Item(1).ColumnWidth = 14

# Leave selection on A4, and this sheet as the active tab
$sel = $newSheet.Range("A4").Select()

Write-Host "Done"
